$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'42.665.35"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.98%  '
$ws.Range('D3').Value = "'2.354.70"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = "'325.29"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.32%  '
$ws.Range('D6').Value = "'100.57"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -8.83%  '
$ws.Range('D7').Value = "'0.636"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.17%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = "'0.624"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.01%  '
$ws.Range('D10').Value = "'39.89"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.89%  '
$ws.Range('D11').Value = "'0.0922"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.90%  '
$ws.Range('D12').Value = "'8.42"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.89%  '
$ws.Range('D13').Value = "'0.999"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.47%  '
$ws.Range('E14').Value = '  -0.01%  '
$ws.Range('D15').Value = "'16.55"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.78%  '
$ws.Range('D16').Value = "'2.710.13"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('D17').Value = "'2.346.55"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.23%  '
$ws.Range('D18').Value = "'8.07"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +11.23%  '
$ws.Range('D19').Value = "'42.622.11"
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Value = "'0.0000107"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.32%  '
$ws.Range('D21').Value = "'76.21"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.73%  '
$ws.Range('D22').Value = "'3.72"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +7.18%  '
$ws.Range('D23').Value = "'266.46"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.31%  '
$ws.Range('D24').Value = "'2.31"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -10.70%  '
$ws.Range('D25').Value = "'10.06"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +10.04%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').Value = "'11.46"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.01%  '
$ws.Range('D28').Value = "'22.97"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.46%  '
$ws.Range('E29').Value = '  -2.18%  '
$ws.Range('D30').Value = "'175.67"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.01%  '
$ws.Range('E31').Value = '  -2.45%  '
$ws.Range('D32').Value = "'0.0900"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.51%  '
$ws.Range('D33').Value = "'35.29"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -10.14%  '
$ws.Range('D34').Value = "'6.02"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.40%  '
$ws.Range('E35').Value = '  -0.29%  '
$ws.Range('D36').Value = "'4.57"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -8.68%  '
$ws.Range('D37').Value = "'0.0357"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.24%  '
$ws.Range('D38').Value = "'2.93"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.10%  '
$ws.Range('D39').Value = "'0.106"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.06%  '
$ws.Range('D40').Value = "'3.79"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -9.41%  '
$ws.Range('D41').Value = "'1.51"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.28%  '
$ws.Range('D42').Value = "'0.235"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('D43').Value = "'70.00"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.72%  '
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('D45').Value = "'119.00"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.73%  '
$ws.Range('D46').Value = "'90.40"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +26.86%  '
$ws.Range('D47').Value = "'11.82"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -8.05%  '
$ws.Range('D48').Value = "'5.50"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.70%  '
$ws.Range('D49').Value = "'9.18"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.38%  '
$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D50').Value = "'1.26"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.19%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = "'1.563.37"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.28%  '
